$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the new "ECs" sending-cluster row (new TPM-derived values).
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Dsg1a"
$ws.Cells.Item(2, 3).Value = "Dsc3"
$ws.Cells.Item(2, 4).Value = "MuSCs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.004083333333333334
$ws.Cells.Item(2, 8).Value = 0.01225
$ws.Cells.Item(2, 9).Value = 0.1210581968752162
$ws.Cells.Item(2, 10).Value = 0.1210581968752162
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.003147333333333334
$ws.Cells.Item(2, 14).Value = 0.009442000000000001
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(2, 17).Value = 0.00001285161111111111
$ws.Cells.Item(2, 18).Value = 0.0001156645
$ws.Cells.Item(2, 19).Value = 0.1210581968752162
$ws.Cells.Item(2, 20).Value = 0.1210581968752162

# Row 3 is the original "MuSCs" sending-cluster row, re-derived with the
# updated edge-specificity values from the new TPM run.
$ws.Cells.Item(3, 1).Value = "MuSCs"
$ws.Cells.Item(3, 2).Value = "Dsg1a"
$ws.Cells.Item(3, 3).Value = "Dsc3"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.029647
$ws.Cells.Item(3, 8).Value = 0.08894100000000001
$ws.Cells.Item(3, 9).Value = 0.8789418031247839
$ws.Cells.Item(3, 10).Value = 0.8789418031247839
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.003147333333333334
$ws.Cells.Item(3, 14).Value = 0.009442000000000001
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 0.00009330899133333336
$ws.Cells.Item(3, 18).Value = 0.0008397809220000002
$ws.Cells.Item(3, 19).Value = 0.8789418031247839
$ws.Cells.Item(3, 20).Value = 0.8789418031247839
